# Apply cryptos list update (row values + 3-row reorder for Maker/Filecoin/OKB)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="58.086.65"'
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Formula = '="3.128.78"'
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("D4").Formula = '="1.00"'
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Formula = '="532.90"'
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  +1.13%  "

$ws.Range("D6").Formula = '="138.25"'
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Formula = '="1.00"'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Formula = '="3.130.61"'
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("D9").Formula = '="0.464"'
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  +4.59%  "

$ws.Range("D10").Formula = '="7.31"'
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = "  +2.41%  "

$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("E12").Value = "  +3.72%  "

$ws.Range("D13").Formula = '="3.664.71"'
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("E14").Value = "  +1.42%  "

$ws.Range("D15").Formula = '="25.66"'
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Formula = '="58.136.85"'
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = "  +0.71%  "

$ws.Range("D18").Formula = '="3.130.11"'
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").Formula = '="12.66"'
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("D21").Formula = '="8.10"'
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  +2.75%  "

$ws.Range("D22").Formula = '="354.18"'
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("D24").Formula = '="69.06"'
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null

$ws.Range("D25").Formula = '="0.504"'
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").Formula = '="0.0₃0880"'
$ws.Range("D28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Value = "  -3.35%  "

$ws.Range("D29").Formula = '="7.28"'
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  -2.24%  "

$ws.Range("D30").Formula = '="6.17"'
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = "  -0.71%  "

$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").Formula = '="21.41"'
$ws.Range("D32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Value = "  +1.34%  "

$ws.Range("D33").Formula = '="5.00"'
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("E34").Value = "  -3.36%  "

$ws.Range("D35").Formula = '="158.69"'
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null

$ws.Range("D36").Formula = '="6.08"'
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("D37").Formula = '="25.80"'
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  -1.05%  "

$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("E39").Value = "  +5.52%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Formula = '="2.457.86"'
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  +6.15%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Formula = '="3.99"'
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  -4.77%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Formula = '="37.60"'
$ws.Range("D44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Value = "  +3.02%  "

$ws.Range("D45").Formula = '="3.170.71"'
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").Formula = '="19.84"'
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  -2.16%  "

$ws.Range("D51").Formula = '="0.740"'
$ws.Range("D51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  -2.80%  "

$excel.CutCopyMode = 0

